# "Le login/register + template de calendrier"
# Fill in the next batch of journal-de-bord entries (rows 43-47 of the
# "Tableau1" log table) and move the active selection/scroll position on
# to where work continues (row 48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Tache (col B) -----------------------------------------------------
# Entered out of row order (matches how the shared-string table grew),
# the empty row at 43 ("finission du login") gets journalled after the
# register/login adjustments below it were already written up.
$ws.Range("B44").Value = "Code du register"
$ws.Range("B45").Value = "Ajustement visuel du register"
$ws.Range("B46").Value = "Ajustement visuel du login"
$ws.Range("E45").Value = "avec une case d'erreur si la connexion echoue"
$ws.Range("B43").Value = "finission du login"
$ws.Range("B47").Value = "Ajout du calendrier trouver en exemple "
$ws.Range("E47").Value = "en plus de tout les ajustement a faire des noms des variables du decorticage du php dans l'html, et de mettre le tout en structure MVC"
$ws.Range("E46").Value = "avec une case d'erreur si la connexion echoue"

# --- Date (col C) / Temps en minutes (col D) ---------------------------
$ws.Range("C43").Value = 44333
$ws.Range("D43").Value = 60

$ws.Range("C44").Value = 44333
$ws.Range("D44").Value = 120

$ws.Range("C45").Value = 44333
$ws.Range("D45").Value = 60

$ws.Range("C46").Value = 44333
$ws.Range("D46").Value = 15

$ws.Range("C47").Value = 44333
$ws.Range("D47").Value = 120

# D46 picks up the wrapped/centered style shared by the rest of column D
# (previously it still carried the older un-wrapped style).
$ws.Range("D46").WrapText = $true

# Row 47's description is long enough to wrap onto two lines.
$ws.Rows.Item(47).RowHeight = 30

# Move the selection/scroll on to the next blank row, ready for the next entry.
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 1
$ws.Range("B48").Select() | Out-Null
